# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the second data row (9da7e663-...) on both the zh-cn and de-de
# worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 02:59:34"
$wsZhCn.Range("H3").Value = "2016-03-17 03:00:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 02:59:41"
$wsDeDe.Range("H3").Value = "2016-03-17 03:00:31"
